$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mechanical section: use accurate unit costs derived from bulk purchase prices ---

# Prongs: $348 for a box of 150 -> per-unit cost
$ws.Range("B15").Formula = "=348/150"

# Zipties: $12.95 for a bag of 50 -> per-unit cost
$ws.Range("B16").Formula = "=12.95/50"

# --- Add "Product Link" hyperlinks for the mechanical parts (column F) ---
# F14: case supplier link
$ws.Range("F14").Value = "http://www.omnicase.com/"
$ws.Hyperlinks.Add($ws.Range("F14"), "http://www.omnicase.com/")
# Re-apply the same look used by the other Product Link cells (e.g. F3) so the
# cell keeps the existing shared "Hyperlink" style instead of a fresh one.
$ws.Range("F3").Copy()
$ws.Range("F14").PasteSpecial(-4122)

# F15: foam-cutting supplier link
$ws.Range("F15").Value = "https://www.bigbluesaw.com/"
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.bigbluesaw.com/")
$ws.Range("F3").Copy()
$ws.Range("F15").PasteSpecial(-4122)

# F16: zipties supplier link
$ws.Range("F16").Value = "http://goo.gl/Xpt3aH"
$ws.Hyperlinks.Add($ws.Range("F16"), "http://goo.gl/Xpt3aH")
$ws.Range("F3").Copy()
$ws.Range("F16").PasteSpecial(-4122)

# F17: nuts & bolts sourced in lab, not purchased online -> plain note, no link
$ws.Range("F17").Value = "in lab"

# Leave the cursor where the author left it after the edit
$ws.Range("F18").Select() | Out-Null
